# "Generate Report for Archive"
#
# The localization-status report is regenerated: every cell whose status
# was "Ready for handoff" is now "In Translation" (Overview!E2:F4,
# zh-cn!C2:C4, de-de!C2:C4 all share that one string), and the
# "Status"/per-language columns that held that text are narrower to fit
# the new (shorter) value.

$wb = $excel.ActiveWorkbook

# --- 1. Status text: "Ready for handoff" -> "In Translation" -----------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Range("E4").Value = "In Translation"
$overview.Range("F4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"
$dede.Range("C4").Value = "In Translation"

# --- 2. Shrink the columns that held that status text to fit it --------
# (was ~17.22 chars wide for "Ready for handoff"; "In Translation" is
# shorter, so the columns narrow to ~13.41 chars wide)
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
